$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "testliz"
$ws.Range("B10").Value = "test"

$ws.Range("A11").Value = "testproft"
$ws.Range("B11").Value = "test"

$ws.Range("A12").Value = "test1"
$ws.Range("B12").Value = "test"
